$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q4" right after "总计".
#    We do this by duplicating the existing "2022-Q3" sheet (which carries all
#    the right sheet-level formatting / column layout) and placing the copy
#    immediately after "总计". All the other quarter sheets shift right by one
#    position automatically, keeping their own (unchanged) data.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$sourceTemplate = $wb.Worksheets.Item("2022-Q3")
$sourceTemplate.Copy($null, $total)

# The copy is placed right after $total, i.e. becomes worksheet #2.
$newQ4 = $wb.Worksheets.Item(2)
$newQ4.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2. Fill in the new "2022-Q4" sheet with the updated fund holdings.
#    Columns D/E/F (and G where noted) are stored as text in this workbook,
#    so force a text number format before assigning numeric-looking strings
#    to keep their literal formatting (e.g. "0.00", leading zeros in codes).
# ---------------------------------------------------------------------------
$textCols = "B","D","E","F"
foreach ($col in $textCols) {
    $newQ4.Range("$col 2:$col 5" -replace ' ', '').NumberFormat = "@"
}

# Row 2 - 003396 / 东方红优享红利混合A
$newQ4.Range("B2").Value = "003396"
$newQ4.Range("C2").Value = "东方红优享红利混合A"
$newQ4.Range("D2").Value = "14.24"
$newQ4.Range("E2").Value = "81.32"
$newQ4.Range("F2").Value = "2.91"
$newQ4.Range("G2").NumberFormat = "@"
$newQ4.Range("G2").Value = "0.4144"
$newQ4.Range("H2").Value = 9

# Row 3 - 001564 / 东方红京东大数据灵活配置混合A
$newQ4.Range("B3").Value = "001564"
$newQ4.Range("C3").Value = "东方红京东大数据灵活配置混合A"
$newQ4.Range("D3").Value = "9.23"
$newQ4.Range("E3").Value = "75.39"
$newQ4.Range("F3").Value = "2.87"
$newQ4.Range("G3").NumberFormat = "@"
$newQ4.Range("G3").Value = "0.2649"
$newQ4.Range("H3").Value = 10

# Row 4 - 017535 / 东方红京东大数据灵活配置混合C (G4 is a real 0, not text)
$newQ4.Range("B4").Value = "017535"
$newQ4.Range("C4").Value = "东方红京东大数据灵活配置混合C"
$newQ4.Range("D4").Value = "0.00"
$newQ4.Range("E4").Value = "75.39"
$newQ4.Range("F4").Value = "2.87"
$newQ4.Range("G4").NumberFormat = "General"
$newQ4.Range("G4").Value = 0
$newQ4.Range("H4").Value = 10

# Row 5 - 017536 / 东方红优享红利混合C (G5 is a real 0, not text)
$newQ4.Range("B5").Value = "017536"
$newQ4.Range("C5").Value = "东方红优享红利混合C"
$newQ4.Range("D5").Value = "0.00"
$newQ4.Range("E5").Value = "81.32"
$newQ4.Range("F5").Value = "2.91"
$newQ4.Range("G5").NumberFormat = "General"
$newQ4.Range("G5").Value = 0
$newQ4.Range("H5").Value = 9

# ---------------------------------------------------------------------------
# 3. Update the "总计" summary sheet: push the existing rows down by one and
#    insert the new 2022-Q4 totals in row 2. Column A is a plain running
#    index (0,1,2,3,...) that is independent of which quarter's data sits in
#    the row, so only columns B:D need to move - A simply grows by one row.
# ---------------------------------------------------------------------------
$total.Range("A5").Copy($total.Range("A6"))
$total.Range("A6").Value = 4

$total.Range("B5:D5").Copy($total.Range("B6:D6"))
$total.Range("B4:D4").Copy($total.Range("B5:D5"))
$total.Range("B3:D3").Copy($total.Range("B4:D4"))
$total.Range("B2:D2").Copy($total.Range("B3:D3"))

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.68

# ---------------------------------------------------------------------------
# 4. Leave "总计" as the active sheet (matches the original workbook view).
# ---------------------------------------------------------------------------
$total.Activate()
